$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 8774833
$ws.Range("I100").Value = 27779370
$ws.Range("J100").Value = 3507.6924
$ws.Range("K100").Value = 27779370
$ws.Range("L100").Value = 3507.6924
$ws.Range("M100").Value = -27778829
$ws.Range("N100").Value = -4589.6924
$ws.Range("H103").Value = 1019.7
$ws.Range("I103").Value = 700
$ws.Range("J103").Value = 1232.8334
$ws.Range("K103").Value = 2100
$ws.Range("L103").Value = 3698.5002
$ws.Range("M103").Value = -1514
$ws.Range("N103").Value = -4870.5002
$ws.Range("H132").Value = 1494.4419
$ws.Range("I132").Value = 1277.9656
$ws.Range("J132").Value = 1942.8572
$ws.Range("K132").Value = 3833.8968
$ws.Range("L132").Value = 5828.571599999999
$ws.Range("M132").Value = -1303.8968
$ws.Range("N132").Value = -10888.5716
$ws.Range("H135").Value = 1448.4286
$ws.Range("I135").Value = 1131.2
$ws.Range("J135").Value = 2746.182
$ws.Range("K135").Value = 10180.8
$ws.Range("L135").Value = 24715.638
$ws.Range("M135").Value = -7645.800000000001
$ws.Range("N135").Value = -29785.638
$ws.Range("H137").Value = 1320.0944
$ws.Range("I137").Value = 1054.5526
$ws.Range("K137").Value = 3163.6578
$ws.Range("M137").Value = -613.6578
$ws.Range("H138").Value = 2244.5889
$ws.Range("I138").Value = 950.70734
$ws.Range("J138").Value = 3327.2246
$ws.Range("K138").Value = 2852.12202
$ws.Range("L138").Value = 9981.6738
$ws.Range("M138").Value = 2287.87798
$ws.Range("N138").Value = -20261.6738
$ws.Range("H141").Value = 1387.3478
$ws.Range("I141").Value = 911.65717
$ws.Range("J141").Value = 2900.9092
$ws.Range("K141").Value = 2734.97151
$ws.Range("L141").Value = 8702.7276
$ws.Range("M141").Value = 2445.02849
$ws.Range("N141").Value = -19062.7276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").Value = ""
$ws.Range("H32").Value = 3474.96
$ws.Range("I32").Value = 2460.5854
$ws.Range("J32").Value = 8096
$ws.Range("K32").Value = 2460.5854
$ws.Range("L32").Value = 8096
$ws.Range("M32").Value = -2173.5854
$ws.Range("N32").Value = -8670
$ws.Range("H61").Value = 3784.8572
$ws.Range("I61").Value = 4902.0713
$ws.Range("J61").Value = 1550.4286
$ws.Range("K61").Value = 4902.0713
$ws.Range("L61").Value = 1550.4286
$ws.Range("M61").Value = -4690.0713
$ws.Range("N61").Value = -1974.4286
$ws.Range("H102").Value = 9260862
$ws.Range("I102").Value = 9260862
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 9260862
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -9259240
$ws.Range("N102").Value = ""
$ws.Range("H110").Value = 690
$ws.Range("I110").Value = 690
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 690
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1355
$ws.Range("N110").Value = ""
$ws.Range("H132").Value = 1641569.5
$ws.Range("I132").Value = 1478.7391
$ws.Range("J132").Value = 6671181
$ws.Range("K132").Value = 4436.2173
$ws.Range("L132").Value = 20013543
$ws.Range("M132").Value = -1906.2173
$ws.Range("N132").Value = -20018603
$ws.Range("H136").Value = 3784.8572
$ws.Range("I136").Value = 4902.0713
$ws.Range("J136").Value = 1550.4286
$ws.Range("K136").Value = 14706.2139
$ws.Range("L136").Value = 4651.2858
$ws.Range("M136").Value = -12156.2139
$ws.Range("N136").Value = -9751.2858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1267.7273
$ws.Range("I94").Value = 454.6154
$ws.Range("K94").Value = 454.6154
$ws.Range("M94").Value = -3.615400000000022
$ws.Range("H99").Value = 58824804
$ws.Range("I99").Value = 76924030
$ws.Range("K99").Value = 76924030
$ws.Range("M99").Value = -76922532
$ws.Range("H134").Value = 4988.5137
$ws.Range("I134").Value = 6832.45
$ws.Range("J134").Value = 2819.1765
$ws.Range("K134").Value = 20497.35
$ws.Range("L134").Value = 8457.529500000001
$ws.Range("M134").Value = -17962.35
$ws.Range("N134").Value = -13527.5295

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 308531.6
$ws.Range("I31").Value = 1803.5483
$ws.Range("J31").Value = 783960.1
$ws.Range("K31").Value = 1803.5483
$ws.Range("L31").Value = 783960.1
$ws.Range("M31").Value = -1508.5483
$ws.Range("N31").Value = -784550.1
$ws.Range("H34").Value = 308531.6
$ws.Range("I34").Value = 1803.5483
$ws.Range("J34").Value = 783960.1
$ws.Range("K34").Value = 1803.5483
$ws.Range("L34").Value = 783960.1
$ws.Range("M34").Value = -1601.5483
$ws.Range("N34").Value = -784364.1
$ws.Range("H58").Value = 1731.4615
$ws.Range("I58").Value = 1079.95
$ws.Range("K58").Value = 1079.95
$ws.Range("M58").Value = -876.95
$ws.Range("H107").Value = 15152116
$ws.Range("I107").Value = 17544464
$ws.Range("J107").Value = 581
$ws.Range("K107").Value = 17544464
$ws.Range("L107").Value = 581
$ws.Range("M107").Value = -17542544
$ws.Range("N107").Value = -4421
$ws.Range("H132").Value = 1937.2858
$ws.Range("I132").Value = 1584.8667
$ws.Range("J132").Value = 2818.3333
$ws.Range("K132").Value = 4754.6001
$ws.Range("L132").Value = 8454.999899999999
$ws.Range("M132").Value = -2224.6001
$ws.Range("N132").Value = -13514.9999
$ws.Range("H134").Value = 1229.129
$ws.Range("I134").Value = 1249.0667
$ws.Range("J134").Value = 1176.3529
$ws.Range("K134").Value = 3747.2001
$ws.Range("L134").Value = 3529.0587
$ws.Range("M134").Value = -1212.2001
$ws.Range("N134").Value = -8599.058700000001
$ws.Range("H136").Value = 1731.4615
$ws.Range("I136").Value = 1079.95
$ws.Range("K136").Value = 3239.85
$ws.Range("M136").Value = -689.8500000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 41667736
$ws.Range("I113").Value = 76923864
$ws.Range("J113").Value = 1403.6364
$ws.Range("K113").Value = 76923864
$ws.Range("L113").Value = 1403.6364
$ws.Range("M113").Value = -76921694
$ws.Range("N113").Value = -5743.6364
$ws.Range("H132").Value = 1862.8654
$ws.Range("I132").Value = 1334.1389
$ws.Range("J132").Value = 3052.5
$ws.Range("K132").Value = 4002.4167
$ws.Range("L132").Value = 9157.5
$ws.Range("M132").Value = -1472.4167
$ws.Range("N132").Value = -14217.5
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8099239.5
$ws.Range("I132").Value = 10690088
$ws.Range("J132").Value = 2838
$ws.Range("K132").Value = 32070264
$ws.Range("L132").Value = 8514
$ws.Range("M132").Value = -32067734
$ws.Range("N132").Value = -13574
$ws.Range("H136").Value = 7336.3267
$ws.Range("I136").Value = 4579.275
$ws.Range("J136").Value = 19589.889
$ws.Range("K136").Value = 13737.825
$ws.Range("L136").Value = 58769.667
$ws.Range("M136").Value = -11187.825
$ws.Range("N136").Value = -63869.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 86668056
$ws.Range("I107").Value = 166668860
$ws.Range("J107").Value = 6667242.5
$ws.Range("K107").Value = 500006580
$ws.Range("L107").Value = 20001727.5
$ws.Range("M107").Value = -500004660
$ws.Range("N107").Value = -20005567.5
$ws.Range("H113").Value = 1275.2084
$ws.Range("I113").Value = 1314.7693
$ws.Range("J113").Value = 1228.4546
$ws.Range("K113").Value = 3944.3079
$ws.Range("L113").Value = 3685.3638
$ws.Range("M113").Value = -1774.3079
$ws.Range("N113").Value = -8025.3638
$ws.Range("H122").Value = 1734.079
$ws.Range("I122").Value = 1710.1786
$ws.Range("J122").Value = 1801
$ws.Range("K122").Value = 5130.5358
$ws.Range("L122").Value = 5403
$ws.Range("M122").Value = -2680.5358
$ws.Range("N122").Value = -10303
$ws.Range("H126").Value = 943.4
$ws.Range("I126").Value = 547.7143
$ws.Range("J126").Value = 1866.6666
$ws.Range("K126").Value = 1643.1429
$ws.Range("L126").Value = 5599.9998
$ws.Range("M126").Value = 826.8571000000002
$ws.Range("N126").Value = -10539.9998
$ws.Range("H132").Value = 12757.976
$ws.Range("I132").Value = 15682.121
$ws.Range("K132").Value = 47046.363
$ws.Range("M132").Value = -44516.363
$ws.Range("H136").Value = 8067097.5
$ws.Range("I136").Value = 2814.4211
$ws.Range("K136").Value = 8443.263300000001
$ws.Range("M136").Value = -5893.263300000001
